# Read excel file and add numbers to list
# The "date" sheet holds nine-ish metrics per period in columns A-E.
# Columns A and B were entered as fractions (e.g. 0.87) and need to be
# re-expressed in the same unit as the rest of the sheet (x1000); column D
# was entered x10 too big and needs /10; column E was entered x1000 too
# big and needs /1000. A few cells keep the arithmetic as a live formula
# (mirroring how the author actually typed them in), the rest are plain
# corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Formula = "=0.87*10^3"
$ws.Range("B2").Formula = "=0.82*10^3"
$ws.Range("D2").Value = 331.538
$ws.Range("E2").Value = 488.898

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Formula = "=0.87*10^3"
$ws.Range("B3").Value = 820
$ws.Range("D3").Value = 343.5
$ws.Range("E3").Value = 568.687

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Formula = "=0.86*10^3"
$ws.Range("B4").Value = 810
$ws.Range("D4").Formula = "=3625.21*10^-1"
$ws.Range("E4").Value = 685.237

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = 830
$ws.Range("B5").Value = 800
$ws.Range("D5").Value = 374.438
$ws.Range("E5").Value = 757.239

# --- Row 6 ---------------------------------------------------------------
$ws.Range("A6").Value = 820
$ws.Range("B6").Value = 790
$ws.Range("D6").Value = 387.743
$ws.Range("E6").Value = 827.08

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = 790
$ws.Range("B7").Value = 780
$ws.Range("D7").Value = 400.399
$ws.Range("E7").Value = 760.056

# --- Row 8 ---------------------------------------------------------------
$ws.Range("A8").Value = 770
$ws.Range("B8").Value = 780
$ws.Range("D8").Value = 415.086
$ws.Range("E8").Value = 620.855

# --- Row 9 ---------------------------------------------------------------
$ws.Range("A9").Value = 550
$ws.Range("B9").Value = 990
$ws.Range("D9").Value = 429.052
$ws.Range("E9").Value = 542.843

# --- Row 10 --------------------------------------------------------------
$ws.Range("A10").Value = 360
$ws.Range("B10").Value = 1160
$ws.Range("D10").Value = 452.789
$ws.Range("E10").Value = 438.404

# --- Row 11 --------------------------------------------------------------
$ws.Range("A11").Value = 340
$ws.Range("B11").Value = 1180
$ws.Range("D11").Value = 483.476
$ws.Range("E11").Value = 383.301

# --- Row 12 --------------------------------------------------------------
$ws.Range("A12").Value = 640
$ws.Range("B12").Value = 870
$ws.Range("D12").Value = 518.163
$ws.Range("E12").Value = 328.88

# Column E now holds values below 1000 (it used to hold big integers), so
# give it two decimal places instead of the old integer format.
$ws.Range("E2:E12").NumberFormat = "0.00"

# Mirror the author's final selection (they'd just dragged over the table).
$ws.Range("A2:E12").Select() | Out-Null
